# This script updates the weekly "Hortaliza, Agrícola del Norte S.A. de Arica -
# Alcachofa" consolidated price sheet. The underlying report re-pulled the
# daily rows in a different order/contents (weekly refresh), so every data
# column (Fecha, Variedad, Calidad, Volumen, Precio mínimo/máximo/promedio,
# Unidad de comercialización, Origen, Precio $/Kg, Kg o Unidades) for rows
# 2-25 is rewritten to its new value. Columns that never change (Mercado ID,
# Mercado, Región, Codreg, Categoría ID, Categoría, Clasificación) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44167
$ws.Range("H2").Value = "Española"
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 13500
$ws.Range("N2").Value = "`$/caja 30 unidades"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 450
$ws.Range("Q2").Value = 30
$ws.Range("D3").Value = 44489
$ws.Range("H3").Value = "Madrigal"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 13500
$ws.Range("N3").Value = "`$/caja 40 unidades"
$ws.Range("P3").Value = 338
$ws.Range("Q3").Value = 40
$ws.Range("D4").Value = 44405
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 21000
$ws.Range("L4").Value = 22000
$ws.Range("M4").Value = 21500
$ws.Range("P4").Value = 538
$ws.Range("D5").Value = 44370
$ws.Range("H5").Value = "Argentina(o)"
$ws.Range("J5").Value = 140
$ws.Range("K5").Value = 20000
$ws.Range("L5").Value = 21000
$ws.Range("M5").Value = 20429
$ws.Range("N5").Value = "`$/caja 50 unidades"
$ws.Range("P5").Value = 409
$ws.Range("Q5").Value = 50
$ws.Range("D6").Value = 44370
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 22000
$ws.Range("L6").Value = 23000
$ws.Range("M6").Value = 22500
$ws.Range("P6").Value = 562
$ws.Range("D7").Value = 44412
$ws.Range("H7").Value = "Symphony"
$ws.Range("J7").Value = 240
$ws.Range("K7").Value = 21000
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = 21500
$ws.Range("P7").Value = 538
$ws.Range("D8").Value = 44391
$ws.Range("H8").Value = "Madrigal"
$ws.Range("J8").Value = 140
$ws.Range("D9").Value = 44363
$ws.Range("D10").Value = 44468
$ws.Range("H10").Value = "Argentina(o)"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 17000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 17500
$ws.Range("N10").Value = "`$/caja 50 unidades"
$ws.Range("P10").Value = 350
$ws.Range("Q10").Value = 50
$ws.Range("D12").Value = 44160
$ws.Range("J12").Value = 160
$ws.Range("D13").Value = 44377
$ws.Range("K13").Value = 20000
$ws.Range("L13").Value = 21000
$ws.Range("M13").Value = 20333
$ws.Range("P13").Value = 508
$ws.Range("D14").Value = 44377
$ws.Range("H14").Value = "Symphony"
$ws.Range("J14").Value = 60
$ws.Range("D15").Value = 44356
$ws.Range("H15").Value = "Argentina(o)"
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 19000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 19500
$ws.Range("N15").Value = "`$/caja 50 unidades"
$ws.Range("O15").Value = "Región de Coquimbo"
$ws.Range("P15").Value = 390
$ws.Range("Q15").Value = 50
$ws.Range("D16").Value = 44435
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 19000
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = 19500
$ws.Range("P16").Value = 488
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 21000
$ws.Range("L17").Value = 22000
$ws.Range("M17").Value = 21500
$ws.Range("N17").Value = "`$/caja 40 unidades"
$ws.Range("P17").Value = 538
$ws.Range("Q17").Value = 40
$ws.Range("H18").Value = "Madrigal"
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = 19000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 19333
$ws.Range("N18").Value = "`$/caja 50 unidades"
$ws.Range("P18").Value = 387
$ws.Range("Q18").Value = 50
$ws.Range("D19").Value = 44384
$ws.Range("H19").Value = "Symphony"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 20000
$ws.Range("L19").Value = 21000
$ws.Range("M19").Value = 20400
$ws.Range("N19").Value = "`$/caja 40 unidades"
$ws.Range("P19").Value = 510
$ws.Range("Q19").Value = 40
$ws.Range("D20").Value = 44706
$ws.Range("J20").Value = 250
$ws.Range("D21").Value = 44482
$ws.Range("J21").Value = 200
$ws.Range("D22").Value = 44398
$ws.Range("J22").Value = 170
$ws.Range("K22").Value = 21000
$ws.Range("L22").Value = 22000
$ws.Range("M22").Value = 21500
$ws.Range("P22").Value = 538
$ws.Range("D23").Value = 44426
$ws.Range("J23").Value = 150
$ws.Range("K23").Value = 19000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 19500
$ws.Range("P23").Value = 488
$ws.Range("D24").Value = 44433
$ws.Range("H24").Value = "Madrigal"
$ws.Range("J24").Value = 160
$ws.Range("K24").Value = 19000
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 19500
$ws.Range("N24").Value = "`$/caja 40 unidades"
$ws.Range("P24").Value = 488
$ws.Range("Q24").Value = 40
$ws.Range("D25").Value = 44483
$ws.Range("J25").Value = 120
